$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts D:K -> F:M)
$ws.Columns("D:E").Insert()

# Copy number formats/styles from the (now-shifted) original D:E columns (now F:G) into new D:E
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new D/E columns with the latest two quarters of data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 102300
$ws.Range("E8").Value = 95800
$ws.Range("D9").Value = 63900
$ws.Range("E9").Value = 58000
$ws.Range("D10").Value = 38400
$ws.Range("E10").Value = 37900
$ws.Range("D12").Value = 13600
$ws.Range("E12").Value = 13900
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 103200
$ws.Range("E17").Value = 100500
$ws.Range("D18").Value = -900
$ws.Range("E18").Value = -4700
$ws.Range("D20").Value = 200
$ws.Range("E20").Value = -400
$ws.Range("D21").Value = 1400
$ws.Range("E21").Value = -3000
$ws.Range("D22").Value = 3500
$ws.Range("E22").Value = 3500
$ws.Range("D23").Value = -4200
$ws.Range("E23").Value = -8700
$ws.Range("D24").Value = 500
$ws.Range("E24").Value = 500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -4600
$ws.Range("E26").Value = -9200
$ws.Range("D27").Value = -4600
$ws.Range("E27").Value = -9200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -200
$ws.Range("E32").Value = 400
$ws.Range("D33").Value = -4600
$ws.Range("E33").Value = -9200
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -4600
$ws.Range("E35").Value = -9200
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 59400
$ws.Range("E41").Value = 68500
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = "NA"
$ws.Range("D43").Value = 86300
$ws.Range("E43").Value = 66400
$ws.Range("D44").Value = 119500
$ws.Range("E44").Value = 117700
$ws.Range("D45").Value = 23900
$ws.Range("E45").Value = 19300
$ws.Range("D46").Value = 289200
$ws.Range("E46").Value = 271900
$ws.Range("D47").Value = 4800
$ws.Range("E47").Value = 4200
$ws.Range("D48").Value = 21100
$ws.Range("E48").Value = 23100
$ws.Range("D49").Value = 58500
$ws.Range("E49").Value = 58600
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 12500
$ws.Range("E52").Value = 11300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 386100
$ws.Range("E54").Value = 369100
$ws.Range("D57").Value = 31400
$ws.Range("E57").Value = 25900
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("D59").Value = 137600
$ws.Range("E59").Value = 134600
$ws.Range("D60").Value = 169000
$ws.Range("E60").Value = 160500
$ws.Range("D61").Value = 136800
$ws.Range("E61").Value = 128900
$ws.Range("D62").Value = 34100
$ws.Range("E62").Value = 32600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 339900
$ws.Range("E66").Value = 322000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -483000
$ws.Range("E72").Value = -478400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 46100
$ws.Range("E76").Value = 47100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -4600
$ws.Range("E81").Value = -9200
$ws.Range("D83").Value = 2000
$ws.Range("E83").Value = 2100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -13800
$ws.Range("E89").Value = -17800
$ws.Range("D91").Value = -600
$ws.Range("E91").Value = -1600
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -600
$ws.Range("E94").Value = -1600
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 8300
$ws.Range("E100").Value = -2400
$ws.Range("D101").Value = 100
$ws.Range("E101").Value = -400
$ws.Range("D102").Value = -6000
$ws.Range("E102").Value = -22200
